$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '29.414.28'
$ws.Cells.Item(2, 4).NumberFormat = 'General'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '  +0.00%  '
$ws.Cells.Item(2, 5).NumberFormat = 'General'

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.849.96'
$ws.Cells.Item(3, 4).NumberFormat = 'General'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '  +0.02%  '
$ws.Cells.Item(3, 5).NumberFormat = 'General'

$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.9997'
$ws.Cells.Item(4, 4).NumberFormat = 'General'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '  +0.09%  '
$ws.Cells.Item(4, 5).NumberFormat = 'General'

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '240.96'
$ws.Cells.Item(5, 4).NumberFormat = 'General'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '  +0.13%  '
$ws.Cells.Item(5, 5).NumberFormat = 'General'

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.6295'
$ws.Cells.Item(6, 4).NumberFormat = 'General'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '  -0.56%  '
$ws.Cells.Item(6, 5).NumberFormat = 'General'

$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '  +0.04%  '
$ws.Cells.Item(7, 5).NumberFormat = 'General'

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.07686'
$ws.Cells.Item(8, 4).NumberFormat = 'General'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '  +1.44%  '
$ws.Cells.Item(8, 5).NumberFormat = 'General'

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.2940'
$ws.Cells.Item(9, 4).NumberFormat = 'General'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '  -0.87%  '
$ws.Cells.Item(9, 5).NumberFormat = 'General'

$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '  +0.14%  '
$ws.Cells.Item(10, 5).NumberFormat = 'General'

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.07748'
$ws.Cells.Item(11, 4).NumberFormat = 'General'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '  +0.47%  '
$ws.Cells.Item(11, 5).NumberFormat = 'General'

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '1.852.98'
$ws.Cells.Item(12, 4).NumberFormat = 'General'
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '  -0.40%  '
$ws.Cells.Item(12, 5).NumberFormat = 'General'

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '5.028'
$ws.Cells.Item(13, 4).NumberFormat = 'General'
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '  +0.45%  '
$ws.Cells.Item(13, 5).NumberFormat = 'General'

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.00001091'
$ws.Cells.Item(14, 4).NumberFormat = 'General'
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '  +8.40%  '
$ws.Cells.Item(14, 5).NumberFormat = 'General'

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.6813'
$ws.Cells.Item(15, 4).NumberFormat = 'General'
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '  -0.81%  '
$ws.Cells.Item(15, 5).NumberFormat = 'General'

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '83.73'
$ws.Cells.Item(16, 4).NumberFormat = 'General'
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '  +0.62%  '
$ws.Cells.Item(16, 5).NumberFormat = 'General'

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '2.103.97'
$ws.Cells.Item(17, 4).NumberFormat = 'General'
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '  -1.20%  '
$ws.Cells.Item(17, 5).NumberFormat = 'General'

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '6.159'
$ws.Cells.Item(18, 4).NumberFormat = 'General'
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '  -0.12%  '
$ws.Cells.Item(18, 5).NumberFormat = 'General'

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '29.443.78'
$ws.Cells.Item(19, 4).NumberFormat = 'General'
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '  -0.01%  '
$ws.Cells.Item(19, 5).NumberFormat = 'General'

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '229.17'
$ws.Cells.Item(20, 4).NumberFormat = 'General'
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '  +0.18%  '
$ws.Cells.Item(20, 5).NumberFormat = 'General'

$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '  -0.12%  '
$ws.Cells.Item(21, 5).NumberFormat = 'General'

$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '  +0.03%  '
$ws.Cells.Item(22, 5).NumberFormat = 'General'

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '7.453'
$ws.Cells.Item(23, 4).NumberFormat = 'General'
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '  -1.71%  '
$ws.Cells.Item(23, 5).NumberFormat = 'General'

$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '  +0.02%  '
$ws.Cells.Item(24, 5).NumberFormat = 'General'

$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '  -0.37%  '
$ws.Cells.Item(25, 5).NumberFormat = 'General'

$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '  -1.13%  '
$ws.Cells.Item(26, 5).NumberFormat = 'General'

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '8.391'
$ws.Cells.Item(27, 4).NumberFormat = 'General'
$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '  +0.00%  '
$ws.Cells.Item(27, 5).NumberFormat = 'General'

$ws.Cells.Item(28, 5).NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '  -0.18%  '
$ws.Cells.Item(28, 5).NumberFormat = 'General'

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '1.310'
$ws.Cells.Item(29, 4).NumberFormat = 'General'
$ws.Cells.Item(29, 5).NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '  +3.99%  '
$ws.Cells.Item(29, 5).NumberFormat = 'General'

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '1.465'
$ws.Cells.Item(30, 4).NumberFormat = 'General'
$ws.Cells.Item(30, 5).NumberFormat = '@'
$ws.Cells.Item(30, 5).Value = '  -0.30%  '
$ws.Cells.Item(30, 5).NumberFormat = 'General'

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.05699'
$ws.Cells.Item(31, 4).NumberFormat = 'General'
$ws.Cells.Item(31, 5).NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '  -0.31%  '
$ws.Cells.Item(31, 5).NumberFormat = 'General'

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '4.115'
$ws.Cells.Item(32, 4).NumberFormat = 'General'

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '4.054'
$ws.Cells.Item(33, 4).NumberFormat = 'General'
$ws.Cells.Item(33, 5).NumberFormat = '@'
$ws.Cells.Item(33, 5).Value = '  +0.60%  '
$ws.Cells.Item(33, 5).NumberFormat = 'General'

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.853'
$ws.Cells.Item(34, 4).NumberFormat = 'General'
$ws.Cells.Item(34, 5).NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '  +0.07%  '
$ws.Cells.Item(34, 5).NumberFormat = 'General'

$ws.Cells.Item(35, 5).NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '  +0.27%  '
$ws.Cells.Item(35, 5).NumberFormat = 'General'

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.7067'
$ws.Cells.Item(36, 4).NumberFormat = 'General'
$ws.Cells.Item(36, 5).NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '  -1.50%  '
$ws.Cells.Item(36, 5).NumberFormat = 'General'

$ws.Cells.Item(37, 5).NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '  -0.02%  '
$ws.Cells.Item(37, 5).NumberFormat = 'General'

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.777'
$ws.Cells.Item(38, 4).NumberFormat = 'General'
$ws.Cells.Item(38, 5).NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '  -0.14%  '
$ws.Cells.Item(38, 5).NumberFormat = 'General'

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.01800'
$ws.Cells.Item(39, 4).NumberFormat = 'General'
$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '  -0.63%  '
$ws.Cells.Item(39, 5).NumberFormat = 'General'

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '1.225.06'
$ws.Cells.Item(40, 4).NumberFormat = 'General'
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '  -3.21%  '
$ws.Cells.Item(40, 5).NumberFormat = 'General'

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '6.443'
$ws.Cells.Item(41, 4).NumberFormat = 'General'
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '  +4.14%  '
$ws.Cells.Item(41, 5).NumberFormat = 'General'

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.9078'
$ws.Cells.Item(42, 4).NumberFormat = 'General'
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '  +0.06%  '
$ws.Cells.Item(42, 5).NumberFormat = 'General'

$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '  +0.02%  '
$ws.Cells.Item(43, 5).NumberFormat = 'General'

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '2.012.78'
$ws.Cells.Item(44, 4).NumberFormat = 'General'
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '  -1.23%  '
$ws.Cells.Item(44, 5).NumberFormat = 'General'

$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '  -0.09%  '
$ws.Cells.Item(45, 5).NumberFormat = 'General'

$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '  -0.23%  '
$ws.Cells.Item(46, 5).NumberFormat = 'General'

$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '  +2.57%  '
$ws.Cells.Item(47, 5).NumberFormat = 'General'

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '7.138'
$ws.Cells.Item(48, 4).NumberFormat = 'General'
$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '  +0.66%  '
$ws.Cells.Item(48, 5).NumberFormat = 'General'

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.4026'
$ws.Cells.Item(49, 4).NumberFormat = 'General'
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '  -0.72%  '
$ws.Cells.Item(49, 5).NumberFormat = 'General'

$ws.Cells.Item(50, 2).Value = 'RenderToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '1.693'
$ws.Cells.Item(50, 4).NumberFormat = 'General'
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '  +0.30%  '
$ws.Cells.Item(50, 5).NumberFormat = 'General'

$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '9.023'
$ws.Cells.Item(51, 4).NumberFormat = 'General'
$ws.Cells.Item(51, 5).NumberFormat = '@'
$ws.Cells.Item(51, 5).Value = '  -1.16%  '
$ws.Cells.Item(51, 5).NumberFormat = 'General'
